# Updates the cryptocurrency price list: refreshed prices/hours for every
# coin, inserted "One" as a new row (shifting TigerCash..BTSEToken down by
# one position), and tweaked a couple of "Best/Worst in 24h" labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking values that are
# stored as text in this workbook. Force a text number format on each such
# cell before assigning its value so Excel doesn't silently convert the
# string to a number (which would also drop significant trailing zeros,
# e.g. "6.500" or "0.00005180").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.87'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '12'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.06'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '12'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.402'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '12'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05935'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '12'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.415'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '12'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.500'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '12'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8114'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '12'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9241'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '12'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1429'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '12'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07431'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '12'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '12'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03073'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '12'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09347'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '12'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.846'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '12'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001566'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '12'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04708'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '12'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005963'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '12'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005918'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '12'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001257'
$ws.Range("E20").Value = '19BitKanKANBestin24h'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '12'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004779'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '12'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008002'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '12'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.563'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '12'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.158'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '12'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3239'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '12'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '12'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002341'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '12'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '12'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '12'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '12'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '12'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '12'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '12'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '12'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '12'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '12'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '12'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '12'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '12'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03905'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '12'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006211'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '12'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '12'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '12'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008322'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '12'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005180'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '12'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '12'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6703'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '12'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002065'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '12'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '12'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '12'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '12'
